$wb = $excel.ActiveWorkbook

# ---- Worksheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H19").Value = 1450.2
$ws.Range("I19").Value = 805
$ws.Range("J19").Value = 1880.3334
$ws.Range("K19").Value = 805
$ws.Range("L19").Value = 1880.3334
$ws.Range("M19").Value = -630
$ws.Range("N19").Value = -2230.3334
$ws.Range("H32").Value = 2724.875
$ws.Range("J32").Value = 2966.5
$ws.Range("L32").Value = 2966.5
$ws.Range("N32").Value = -3618.5
$ws.Range("H98").Value = 1132.6666
$ws.Range("I98").Value = 1132.6666
$ws.Range("K98").Value = 1132.6666
$ws.Range("M98").Value = 365.3334
$ws.Range("H122").Value = 1132.6666
$ws.Range("I122").Value = 1132.6666
$ws.Range("K122").Value = 3397.9998
$ws.Range("M122").Value = -947.9998000000001
$ws.Range("H127").Value = 2522.8333
$ws.Range("J127").Value = 1599
$ws.Range("L127").Value = 4797
$ws.Range("N127").Value = -14717
$ws.Range("H135").Value = 366.30768
$ws.Range("I135").Value = 366.30768
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 3296.76912
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -761.7691199999999
$ws.Range("N135").ClearContents()

# ---- Worksheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H61").Value = 2380
$ws.Range("I61").Value = 2380
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2380
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2168
$ws.Range("N61").ClearContents()
$ws.Range("H122").Value = 2935.4285
$ws.Range("I122").Value = 2189.6
$ws.Range("J122").Value = 4800
$ws.Range("K122").Value = 6568.799999999999
$ws.Range("L122").Value = 14400
$ws.Range("M122").Value = -4118.799999999999
$ws.Range("N122").Value = -19300
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 2582.5386
$ws.Range("I132").Value = 2622.75
$ws.Range("J132").Value = 2100
$ws.Range("K132").Value = 7868.25
$ws.Range("L132").Value = 6300
$ws.Range("M132").Value = -5338.25
$ws.Range("N132").Value = -11360
$ws.Range("H136").Value = 2380
$ws.Range("I136").Value = 2380
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 7140
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -4590
$ws.Range("N136").ClearContents()

# ---- Worksheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H107").Value = 3399.6
$ws.Range("I107").Value = 3399.6
$ws.Range("K107").Value = 3399.6
$ws.Range("M107").Value = -1479.6

# ---- Worksheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H12").Value = 4999.5
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 4999.5
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 4999.5
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -5339.5
$ws.Range("H16").Value = 2865.6667
$ws.Range("I16").Value = 3118.8
$ws.Range("J16").Value = 1600
$ws.Range("K16").Value = 3118.8
$ws.Range("L16").Value = 1600
$ws.Range("M16").Value = -2831.8
$ws.Range("N16").Value = -2174
$ws.Range("H81").Value = 59989.5
$ws.Range("I81").Value = 39980
$ws.Range("J81").Value = 79999
$ws.Range("K81").Value = 39980
$ws.Range("L81").Value = 79999
$ws.Range("M81").Value = -38982
$ws.Range("N81").Value = -81995
$ws.Range("H84").Value = 59989.5
$ws.Range("I84").Value = 39980
$ws.Range("J84").Value = 79999
$ws.Range("K84").Value = 119940
$ws.Range("L84").Value = 239997
$ws.Range("M84").Value = -114948
$ws.Range("N84").Value = -249981
$ws.Range("H94").Value = 5272.778
$ws.Range("I94").Value = 6749.3335
$ws.Range("J94").Value = 4534.5
$ws.Range("K94").Value = 6749.3335
$ws.Range("L94").Value = 4534.5
$ws.Range("M94").Value = -6298.3335
$ws.Range("N94").Value = -5436.5
$ws.Range("H113").Value = 2865.6667
$ws.Range("I113").Value = 3118.8
$ws.Range("J113").Value = 1600
$ws.Range("K113").Value = 3118.8
$ws.Range("L113").Value = 1600
$ws.Range("M113").Value = -948.8000000000002
$ws.Range("N113").Value = -5940
$ws.Range("H133").Value = 25296
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H141").Value = 166995.38
$ws.Range("J141").Value = 166995.38
$ws.Range("L141").Value = 166995.38
$ws.Range("N141").Value = -177355.38

# ---- Worksheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H12").Value = 238.93333
$ws.Range("I12").Value = 242.5
$ws.Range("J12").Value = 231.8
$ws.Range("K12").Value = 727.5
$ws.Range("L12").Value = 695.4000000000001
$ws.Range("M12").Value = -554.5
$ws.Range("N12").Value = -1041.4
$ws.Range("H81").Value = 2266.5
$ws.Range("I81").Value = 1999.5
$ws.Range("J81").Value = 2400
$ws.Range("K81").Value = 5998.5
$ws.Range("L81").Value = 7200
$ws.Range("M81").Value = -4875.5
$ws.Range("N81").Value = -9446
$ws.Range("H84").Value = 2266.5
$ws.Range("I84").Value = 1999.5
$ws.Range("J84").Value = 2400
$ws.Range("K84").Value = 17995.5
$ws.Range("L84").Value = 21600
$ws.Range("M84").Value = -12379.5
$ws.Range("N84").Value = -32832

# ---- Worksheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H7").Value = 4630.375
$ws.Range("I7").Value = 2409
$ws.Range("J7").Value = 8332.666999999999
$ws.Range("K7").Value = 2409
$ws.Range("L7").Value = 8332.666999999999
$ws.Range("M7").Value = -2297
$ws.Range("N7").Value = -8556.666999999999
$ws.Range("H46").Value = 3316.5
$ws.Range("I46").Value = 2954.6
$ws.Range("J46").Value = 3678.4
$ws.Range("K46").Value = 2954.6
$ws.Range("L46").Value = 3678.4
$ws.Range("M46").Value = -2766.6
$ws.Range("N46").Value = -4054.4
$ws.Range("H126").Value = 4630.375
$ws.Range("I126").Value = 2409
$ws.Range("J126").Value = 8332.666999999999
$ws.Range("K126").Value = 7227
$ws.Range("L126").Value = 24998.001
$ws.Range("M126").Value = -4757
$ws.Range("N126").Value = -29938.001
$ws.Range("H132").Value = 7457.6665
$ws.Range("I132").Value = 7681
$ws.Range("K132").Value = 23043
$ws.Range("M132").Value = -20513
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H136").Value = 2920.1428
$ws.Range("I136").Value = 1483
$ws.Range("K136").Value = 4449
$ws.Range("M136").Value = -1899

# ---- Worksheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H18").Value = 20000
$ws.Range("J18").Value = 20000
$ws.Range("L18").Value = 20000
$ws.Range("N18").Value = -20346
$ws.Range("H41").Value = 19979.25
$ws.Range("I41").Value = 19978
$ws.Range("K41").Value = 19978
$ws.Range("M41").Value = -19588
$ws.Range("H122").Value = 1982
$ws.Range("I122").Value = 1937.3334
$ws.Range("J122").Value = 2250
$ws.Range("K122").Value = 5812.0002
$ws.Range("L122").Value = 6750
$ws.Range("M122").Value = -3362.0002
$ws.Range("N122").Value = -11650
